$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final parent outcome measurements added: fill column D (Post Treatment)
# for rows 2-15 with "Not worse", matching the other phase columns.
$ws.Range("D2:D15").Value = "Not worse"

# Move the selection cursor to match the post-edit state
$ws.Range("C19").Select()
